# Load and preprocess IJC (7/17 release)
#
# Fills in the just-released "Actual" / "Revised" initial-jobless-claims
# figures (previously placeholder "nan" values), refreshes the summary
# statistics that shift once the actual prints are known, and updates the
# Bloomberg economist panel with the latest revision from Herrmann
# Forecasting LLC (whose "As of" date moved to 7/09 and whose rank-7
# Citigroup estimate note drops off the list).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headline release figures -------------------------------------------
# Actual (row 7) and Revised (row 9) were "nan" placeholders; now reported.
$ws.Range("C7").Value = 227
$ws.Range("C9").Value = 232

# --- Summary statistics, recomputed now that Actual/Revised are known ---
$ws.Range("C13").Value = "235.20k"   # Average Estimate
$ws.Range("C17").Value = 34          # Qualified Economists
$ws.Range("C18").Value = "4.96k"     # Standard Deviation

# --- Economist estimate table updates -----------------------------------
# Citigroup Inc (row 27) no longer carries a "7th" rank annotation.
$ws.Range("F27").Value = ""

# Herrmann Forecasting LLC (row 58): revised estimate/date, rank unchanged.
$ws.Range("D58").Value = 233
$ws.Range("E58").Value = "7/09/2025"
